$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha Ponto")
$ws.Range("B18").Value = "test"
